$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.405.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "'2.610.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'537.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.08%  "
$ws.Range("D6").Value = "'140.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "'2.617.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D12").Value = "'0.337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "'3.070.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "'59.325.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'20.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "'2.621.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").Value = "'342.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "'4.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'67.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "'7.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").Value = "'0.0₃0746"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.94%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +5.71%  "
$ws.Range("D32").Value = "'5.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "'18.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "'148.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "'4.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "'36.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").Value = "'0.840"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").Value = "'0.839"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'276.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").Value = "'0.599"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("D45").Value = "'0.0962"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("D46").Value = "'10.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'0.0525"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").Value = "'1.952.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D50").Value = "'18.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").Value = "'4.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.25%  "
